$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New trade rows appended to the bottom of the data.
$data = @(
    @(9,  9777.59, 9794.24, 78.05, 77.92, $false, -0.17, 42613.765451388892, $false),
    @(10, 9857.77, 9777.59, 77.739999999999995, 78.38, $false, 0.82, 42614.67260416667, $true),
    @(11, 9739.48, 9857.77, 78.36, 77.42, $false, -1.2, 42615.750023148146, $false)
)

# Reference cell whose formatting (date style) should be copied into the
# new "Date" column cells for each appended row.
$dateFormatSource = $ws.Cells.Item(3, 7)

foreach ($row in $data) {
    $r = $row[0]

    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
    $ws.Cells.Item($r, 5).Value = $row[5]
    $ws.Cells.Item($r, 6).Value = $row[6]

    $dateCell = $ws.Cells.Item($r, 7)
    $dateCell.Value = $row[7]
    $dateFormatSource.Copy()
    $dateCell.PasteSpecial(-4122)

    $ws.Cells.Item($r, 8).Value = $row[8]
}

$excel.CutCopyMode = $false
